$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: the review went from "yes" to "no" ---
$ws.Range("G12").Value = "no"

# --- Row 14: the review went from "yes" to "no" ---
$ws.Range("G14").Value = "no"

# --- Row 15: drop the wrap-text formatting on C15:D15 (copy format from C14:D14,
#     which already has the correct non-wrapping style) ---
$ws.Range("C14:D14").Copy()
$ws.Range("C15:D15").PasteSpecial(-4122)

# --- Rows 16-18: new review rows appended to the table. Copy the row-15
#     formatting down (keeps column styles: A -> s1, C/D -> s2, others default)
#     and then fill in the values. ---
$ws.Range("A15:G15").Copy()
$ws.Range("A16:G18").PasteSpecial(-4122)

$ws.Range("A16").Value = "com.hamxa.shaynachim"
$ws.Range("B16").Value = "bitcoin"
$ws.Range("C16").Value = "oamitay16@gmail.com"
$ws.Range("D16").Value = "efiamid9@gmail.com"
$ws.Range("E16").Value = "27/5/2019 15:59"
$ws.Range("F16").Value = "come as you are"
$ws.Range("G16").Value = "yes"

$ws.Range("A17").Value = "com.hamxa.shaynachim"
$ws.Range("B17").Value = "bitcoin"
$ws.Range("C17").Value = "oferrotberg4@gmail.com"
$ws.Range("D17").Value = "oamitay16@gmail.com"
$ws.Range("E17").Value = "27/5/2019 15:59"
$ws.Range("F17").Value = "This app is my best friend"
$ws.Range("G17").Value = "yes"

$ws.Range("A18").Value = "com.hamxa.shaynachim"
$ws.Range("B18").Value = "bitcoin"
$ws.Range("C18").Value = "segevhod6@gmail.com"
$ws.Range("D18").Value = "oferrotberg4@gmail.com"
$ws.Range("E18").Value = "27/5/2019 15:59"
$ws.Range("F18").Value = "Take me to blockchain ibiza"
$ws.Range("G18").Value = "yes"

# --- Move the active selection to F19, matching the saved cursor position ---
[void]$ws.Range("F19").Select()
